$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.611.56'
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").Value = '1.593.54'
$ws.Range("E3").Value = '  +0.16%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.03%  '

$ws.Range("E7").Value = '  +0.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.245'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0840'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.56%  '

$ws.Range("D12").Value = '1.818.27'
$ws.Range("E12").Value = '  +0.31%  '

$ws.Range("D13").Value = '1.583.87'
$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("E15").Value = '  -1.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").Value = '26.602.73'
$ws.Range("E17").Value = '  -0.18%  '

$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '207.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("E23").Value = '  -3.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.02%  '

$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.115'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("E32").Value = '  -0.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.653'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").Value = '1.279.16'
$ws.Range("E35").Value = '  -3.61%  '

$ws.Range("E36").Value = '  +1.70%  '

$ws.Range("E37").Value = '  -0.91%  '

$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.838'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.28%  '

$ws.Range("E40").Value = '  +0.21%  '

$ws.Range("E41").Value = '  +0.76%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.44%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.787'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("E45").Value = '  +9.49%  '

$ws.Range("D46").Value = '1.731.25'
$ws.Range("E46").Value = '  +0.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.74'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.21%  '

$ws.Range("E48").Value = '  -0.39%  '

$ws.Range("D49").Value = '0.0₆0104'
$ws.Range("E49").Value = '  -1.62%  '

$ws.Range("E50").Value = '  +3.28%  '

$ws.Range("E51").Value = '  -1.22%  '
